$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1008
$ws.Range("K3").Value = 963
$ws.Range("K4").Value = 211
$ws.Range("K6").Value = 1318
$ws.Range("K7").Value = 3556

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 104
$ws.Range("K8").Value = 208
$ws.Range("K9").Value = 18
$ws.Range("K15").Value = 21
$ws.Range("K19").Value = 96
$ws.Range("K20").Value = 86
$ws.Range("K32").Value = 7
$ws.Range("K33").Value = 148
$ws.Range("K37").Value = 114
$ws.Range("K42").Value = 117
$ws.Range("K44").Value = 35
$ws.Range("K45").Value = 4
$ws.Range("K47").Value = 26
$ws.Range("K49").Value = 23
$ws.Range("K52").Value = 87
$ws.Range("K53").Value = 42
$ws.Range("K54").Value = 66
$ws.Range("K57").Value = 7
$ws.Range("K60").Value = 25
$ws.Range("K61").Value = 6
$ws.Range("K63").Value = 14
$ws.Range("K65").Value = 99
$ws.Range("K69").Value = 13
$ws.Range("K70").Value = 9
$ws.Range("K72").Value = 14
$ws.Range("K73").Value = 40
$ws.Range("K78").Value = 54
$ws.Range("K79").Value = 102
$ws.Range("K80").Value = 13
$ws.Range("K88").Value = 50
$ws.Range("K89").Value = 55
$ws.Range("K91").Value = 40
$ws.Range("K96").Value = 58
$ws.Range("K100").Value = 6
$ws.Range("K101").Value = 3556

# West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 8
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 58

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 104

# Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 55

# Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 22
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 87

# Norwood Park
$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K3").Value = 6
$ws.Range("K6").Value = 13

# Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 42

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 208

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 37
$ws.Range("K3").Value = 57
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 148

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 114

# New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 99

# Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 23

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 66

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 69

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 96

# Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 35

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 117

# Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 54

# Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 40

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 38
$ws.Range("K7").Value = 102

# Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 26
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 86

# Wrigleyville
$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 6

# Kenwood
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 26

# Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K3").Value = 1
$ws.Range("K7").Value = 21

# Avalon Park
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 18

# Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K5").Value = 17
$ws.Range("K6").Value = 40

# O'Hare
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("J4").Value = 3
$ws.Range("J6").Value = 9

# United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 11
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 50

# Galewood
$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("K3").Value = 2
$ws.Range("K6").Value = 7

# Mckinley Park
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 7

# Morgan Park
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 25

# Old Town
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K5").Value = 8
$ws.Range("K6").Value = 14

# Jackson Park
$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("K3").Value = 1
$ws.Range("K7").Value = 4

# Rush & Division
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 13

# Mount Greenwood
$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("K2").Value = 3
$ws.Range("K6").Value = 6
